$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: only date changes 44295 -> 44294
$ws.Range("D5").Value = 44294

# Row 6: only date changes 44295 -> 44294
$ws.Range("D6").Value = 44294

# Row 7: date 44294 -> 44316, quality Primera -> Especial,
# volume stays 20, prices 225000/230000/227500 -> 255000/260000/257500,
# origin Región Metropolitana -> Región de O'Higgins, price/kg 506 -> 572
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 255000
$ws.Range("O7").Value = 260000
$ws.Range("P7").Value = 257500
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 572

# Row 8: date 44294 -> 44316, quality Segunda -> Primera,
# volume 16 -> 20, prices 195000/200000/197500 -> 225000/230000/227500,
# origin Región Metropolitana -> Región de O'Higgins, price/kg 439 -> 506
$ws.Range("D8").Value = 44316
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 225000
$ws.Range("O8").Value = 230000
$ws.Range("P8").Value = 227500
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 506

# Row 9: date 44316 -> 44295, quality Especial -> Primera,
# volume stays 20, prices 255000/260000/257500 -> 225000/230000/227500,
# origin Región de O'Higgins -> Región Metropolitana, price/kg 572 -> 506
$ws.Range("D9").Value = 44295
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 225000
$ws.Range("O9").Value = 230000
$ws.Range("P9").Value = 227500
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 506

# Row 10: date 44316 -> 44295, quality Primera -> Segunda,
# volume 20 -> 16, prices 225000/230000/227500 -> 195000/200000/197500,
# origin Región de O'Higgins -> Región Metropolitana, price/kg 506 -> 439
$ws.Range("D10").Value = 44295
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 16
$ws.Range("N10").Value = 195000
$ws.Range("O10").Value = 200000
$ws.Range("P10").Value = 197500
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 439
